# Update "想去人数" (column F) counts on the 展览 and 全部类型 sheets
# to reflect newly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - row => new value for column F
$exhibitionUpdates = @{
    3  = 376
    4  = 440
    5  = 16
    6  = 10
    7  = 248
    8  = 13732
    9  = 74
    10 = 67
    11 = 5555
    12 = 573
    13 = 37
    14 = 32
    17 = 67
    18 = 161
    19 = 743
    20 = 2902
    21 = 19
    22 = 9285
    23 = 1180
    24 = 15
    25 = 24
    26 = 3691
    27 = 235
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" (All Types) - row => new value for column F
$allTypesUpdates = @{
    3  = 376
    5  = 440
    6  = 16
    7  = 10
    8  = 248
    9  = 13732
    10 = 74
    11 = 67
    12 = 5555
    13 = 573
    14 = 37
    15 = 32
    18 = 67
    19 = 161
    20 = 743
    21 = 2902
    22 = 19
    24 = 9285
    25 = 1180
    26 = 15
    27 = 24
    28 = 3691
    29 = 235
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
